$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.544.23'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.452.96'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '508.45'
$ws.Range('E5').Value = '  -2.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.26'
$ws.Range('E6').Value = '  +2.23%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.449.32'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.150'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.63'
$ws.Range('E13').Value = '  -6.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.879.20'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.546.15'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.02'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.432.95'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.35'
$ws.Range('E19').Value = '  +1.09%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '314.91'
$ws.Range('E21').Value = '  -0.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.39'
$ws.Range('E22').Value = '  +4.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  -2.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.42'
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.994'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.537.14'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.62'
$ws.Range('E30').Value = '  +6.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.82'
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0736'
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.20'
$ws.Range('E34').Value = '  +1.34%  '
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.995'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.99'
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.25'
$ws.Range('E39').Value = '  +5.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.87'
$ws.Range('E40').Value = '  +2.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.822'
$ws.Range('E41').Value = '  +4.02%  '
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('E43').Value = '  +1.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '134.64'
$ws.Range('E44').Value = '  +10.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.41'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.01'
$ws.Range('E46').Value = '  +4.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '257.13'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.570'
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0918'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('E50').Value = '  +0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0214'
$ws.Range('E51').Value = '  +1.88%  '
